$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45190 to 45192
# for every data row (rows 2 through 119).
$ws.Range("C2:C119").Value = 45192
